$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 26

$ws.Range("A$newRow").Value = "17-10-2025"
$ws.Range("B$newRow").Value = "The price of gold in India today is ₹13,277 per gram for 24 karat gold, ₹12,170 per gram for 22 karat gold and ₹9,958 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A$newRow" + ":B$newRow").Borders.LineStyle = 1
$ws.Range("B$newRow").WrapText = $true
